# Daily attendance processing - reorder the "Recorded By" (column G) list
# for every data row: reverse the order of the comma-separated recorder
# names/emails (e.g. "a, System" -> "System, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ", "
    $n = $parts.Length

    if ($n -le 1) { continue }

    $reversed = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value2 = [string]::Join(", ", $reversed)
}
